$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2396.8572
$ws.Range("I9").Value = 224.75
$ws.Range("K9").Value = 224.75
$ws.Range("M9").Value = -55.75
$ws.Range("H33").Value = 3882877.2
$ws.Range("I33").Value = 8250234
$ws.Range("J33").Value = 782.6667
$ws.Range("K33").Value = 8250234
$ws.Range("L33").Value = 782.6667
$ws.Range("M33").Value = -8250005
$ws.Range("N33").Value = -1240.6667
$ws.Range("H43").Value = 3499.75
$ws.Range("J43").Value = 3499.75
$ws.Range("L43").Value = 3499.75
$ws.Range("N43").Value = -3637.75
$ws.Range("H51").Value = 55600.2
$ws.Range("I51").Value = 89999.336
$ws.Range("J51").Value = 4001.5
$ws.Range("K51").Value = 89999.336
$ws.Range("L51").Value = 4001.5
$ws.Range("M51").Value = -89515.336
$ws.Range("N51").Value = -4969.5
$ws.Range("H55").Value = 170.06667
$ws.Range("J55").Value = 139.63637
$ws.Range("L55").Value = 139.63637
$ws.Range("N55").Value = -567.6363699999999
$ws.Range("H64").Value = 4970.6875
$ws.Range("I64").Value = 4439.25
$ws.Range("J64").Value = 5502.125
$ws.Range("K64").Value = 4439.25
$ws.Range("L64").Value = 5502.125
$ws.Range("M64").Value = -4191.25
$ws.Range("N64").Value = -5998.125
$ws.Range("H67").Value = 4970.6875
$ws.Range("I67").Value = 4439.25
$ws.Range("J67").Value = 5502.125
$ws.Range("K67").Value = 4439.25
$ws.Range("L67").Value = 5502.125
$ws.Range("M67").Value = -3581.25
$ws.Range("N67").Value = -7218.125
$ws.Range("H74").Value = 5997.5625
$ws.Range("I74").Value = 5639.7856
$ws.Range("K74").Value = 5639.7856
$ws.Range("M74").Value = -4703.7856
$ws.Range("H76").Value = 4925.375
$ws.Range("I76").Value = 4925.375
$ws.Range("K76").Value = 4925.375
$ws.Range("M76").Value = -4610.375
$ws.Range("H77").Value = 5997.5625
$ws.Range("I77").Value = 5639.7856
$ws.Range("K77").Value = 28198.928
$ws.Range("M77").Value = -23518.928
$ws.Range("H79").Value = 4925.375
$ws.Range("I79").Value = 4925.375
$ws.Range("K79").Value = 4925.375
$ws.Range("M79").Value = -3833.375
$ws.Range("H80").Value = 4520.8
$ws.Range("I80").Value = 800
$ws.Range("J80").Value = 5451
$ws.Range("K80").Value = 2400
$ws.Range("L80").Value = 16353
$ws.Range("M80").Value = -1402
$ws.Range("N80").Value = -18349
$ws.Range("H83").Value = 4520.8
$ws.Range("I83").Value = 800
$ws.Range("J83").Value = 5451
$ws.Range("K83").Value = 7200
$ws.Range("L83").Value = 49059
$ws.Range("M83").Value = -2208
$ws.Range("N83").Value = -59043
$ws.Range("H88").Value = 522.6923
$ws.Range("J88").Value = 513.7778
$ws.Range("L88").Value = 513.7778
$ws.Range("N88").Value = -1325.7778
$ws.Range("H91").Value = 522.6923
$ws.Range("J91").Value = 513.7778
$ws.Range("L91").Value = 513.7778
$ws.Range("N91").Value = -3321.7778
$ws.Range("H94").Value = 911.25
$ws.Range("I94").Value = 916.6667
$ws.Range("J94").Value = 895
$ws.Range("K94").Value = 916.6667
$ws.Range("L94").Value = 895
$ws.Range("M94").Value = -465.6667
$ws.Range("N94").Value = -1797
$ws.Range("H99").Value = 1823.2222
$ws.Range("I99").Value = 1542.5714
$ws.Range("J99").Value = 2805.5
$ws.Range("K99").Value = 4627.7142
$ws.Range("L99").Value = 8416.5
$ws.Range("M99").Value = -3129.7142
$ws.Range("N99").Value = -11412.5
$ws.Range("H101").Value = 1794.909
$ws.Range("I101").Value = 621.1429000000001
$ws.Range("J101").Value = 3849
$ws.Range("K101").Value = 1863.4287
$ws.Range("L101").Value = 11547
$ws.Range("M101").Value = -241.4287000000002
$ws.Range("N101").Value = -14791
$ws.Range("H116").Value = 3325.2222
$ws.Range("I116").Value = 3280.25
$ws.Range("K116").Value = 3280.25
$ws.Range("M116").Value = 161.75
$ws.Range("H132").Value = 1359.4412
$ws.Range("I132").Value = 876.8889
$ws.Range("J132").Value = 3220.7144
$ws.Range("K132").Value = 2630.6667
$ws.Range("L132").Value = 9662.143199999999
$ws.Range("M132").Value = -100.6667000000002
$ws.Range("N132").Value = -14722.1432
$ws.Range("H133").Value = 99984.5
$ws.Range("J133").Value = 99984.5
$ws.Range("L133").Value = 99984.5
$ws.Range("N133").Value = -110104.5
$ws.Range("H137").Value = 4217.1763
$ws.Range("I137").Value = 1391.25
$ws.Range("J137").Value = 10999.4
$ws.Range("K137").Value = 4173.75
$ws.Range("L137").Value = 32998.2
$ws.Range("M137").Value = -1623.75
$ws.Range("N137").Value = -38098.2
$ws.Range("H138").Value = 2816.6458
$ws.Range("I138").Value = 1940.1
$ws.Range("K138").Value = 5820.299999999999
$ws.Range("M138").Value = -680.2999999999993

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 205.54546
$ws.Range("I5").Value = 205.54546
$ws.Range("K5").Value = 205.54546
$ws.Range("M5").Value = -93.54545999999999
$ws.Range("H32").Value = 2187.9744
$ws.Range("I32").Value = 2203.8333
$ws.Range("J32").Value = 1997.6666
$ws.Range("K32").Value = 2203.8333
$ws.Range("L32").Value = 1997.6666
$ws.Range("M32").Value = -1916.8333
$ws.Range("N32").Value = -2571.6666
$ws.Range("H61").Value = 4645.381
$ws.Range("I61").Value = 4511
$ws.Range("J61").Value = 4712.5713
$ws.Range("K61").Value = 4511
$ws.Range("L61").Value = 4712.5713
$ws.Range("M61").Value = -4299
$ws.Range("N61").Value = -5136.5713
$ws.Range("H74").Value = 2014589.8
$ws.Range("I74").Value = 1158944.5
$ws.Range("K74").Value = 1158944.5
$ws.Range("M74").Value = -1158070.5
$ws.Range("H77").Value = 2014589.8
$ws.Range("I77").Value = 1158944.5
$ws.Range("K77").Value = 5794722.5
$ws.Range("M77").Value = -5790354.5
$ws.Range("H97").Value = 1058.7693
$ws.Range("I97").Value = 887.63635
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 887.63635
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -391.63635
$ws.Range("N97").Value = -2992
$ws.Range("H117").Value = 58333
$ws.Range("J117").Value = 58333
$ws.Range("L117").Value = 58333
$ws.Range("N117").Value = -67511
$ws.Range("H122").Value = 2978.25
$ws.Range("I122").Value = 2641
$ws.Range("K122").Value = 7923
$ws.Range("M122").Value = -5473
$ws.Range("H132").Value = 100000824
$ws.Range("I132").Value = 1030.25
$ws.Range("K132").Value = 3090.75
$ws.Range("M132").Value = -560.75
$ws.Range("H136").Value = 4645.381
$ws.Range("I136").Value = 4511
$ws.Range("J136").Value = 4712.5713
$ws.Range("K136").Value = 13533
$ws.Range("L136").Value = 14137.7139
$ws.Range("M136").Value = -10983
$ws.Range("N136").Value = -19237.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 205.54546
$ws.Range("I4").Value = 205.54546
$ws.Range("K4").Value = 205.54546
$ws.Range("M4").Value = -90.54545999999999
$ws.Range("H86").Value = 2836.1667
$ws.Range("I86").Value = 2377
$ws.Range("K86").Value = 2377
$ws.Range("M86").Value = -1254
$ws.Range("H89").Value = 2836.1667
$ws.Range("I89").Value = 2377
$ws.Range("K89").Value = 11885
$ws.Range("M89").Value = -6269
$ws.Range("H132").Value = 85780
$ws.Range("J132").Value = 85780
$ws.Range("L132").Value = 85780
$ws.Range("N132").Value = -95900
$ws.Range("H134").Value = 17680698
$ws.Range("I134").Value = 8068623
$ws.Range("K134").Value = 24205869
$ws.Range("M134").Value = -24203334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 5264264
$ws.Range("I19").Value = 7692450
$ws.Range("J19").Value = 3193.5
$ws.Range("K19").Value = 7692450
$ws.Range("L19").Value = 3193.5
$ws.Range("M19").Value = -7692280
$ws.Range("N19").Value = -3533.5
$ws.Range("H24").Value = 5264264
$ws.Range("I24").Value = 7692450
$ws.Range("J24").Value = 3193.5
$ws.Range("K24").Value = 7692450
$ws.Range("L24").Value = 3193.5
$ws.Range("M24").Value = -7692280
$ws.Range("N24").Value = -3533.5
$ws.Range("H31").Value = 1961.6666
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 1961.6666
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H95").Value = 27500
$ws.Range("J95").Value = 27500
$ws.Range("L95").Value = 27500
$ws.Range("N95").Value = -32992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J2").Value = 62
$ws.Range("L2").Value = 372
$ws.Range("N2").Value = -598
$ws.Range("H4").Value = 12625948
$ws.Range("J4").Value = 251237.38
$ws.Range("L4").Value = 753712.14
$ws.Range("N4").Value = -753936.14
$ws.Range("H103").Value = 209.15384
$ws.Range("I103").Value = 172.16667
$ws.Range("J103").Value = 240.85715
$ws.Range("K103").Value = 516.50001
$ws.Range("L103").Value = 722.5714499999999
$ws.Range("M103").Value = 362.49999
$ws.Range("N103").Value = -2480.57145
$ws.Range("H114").Value = 1560
$ws.Range("I114").Value = 436.6
$ws.Range("J114").Value = 2028.0834
$ws.Range("K114").Value = 1309.8
$ws.Range("L114").Value = 6084.2502
$ws.Range("M114").Value = 1944.2
$ws.Range("N114").Value = -12592.2502
$ws.Range("H117").Value = 3972.7
$ws.Range("J117").Value = 4856.857
$ws.Range("L117").Value = 14570.571
$ws.Range("N117").Value = -21454.571
$ws.Range("H121").Value = 1366736
$ws.Range("I121").Value = 2161.2856
$ws.Range("J121").Value = 2731310.8
$ws.Range("K121").Value = 6483.8568
$ws.Range("L121").Value = 8193932.399999999
$ws.Range("M121").Value = -5173.8568
$ws.Range("N121").Value = -8196552.399999999
$ws.Range("H131").Value = 628232.75
$ws.Range("J131").Value = 1400221.9
$ws.Range("L131").Value = 4200665.699999999
$ws.Range("N131").Value = -4210745.699999999
$ws.Range("H134").Value = 10333.8
$ws.Range("I134").Value = 12186.333
$ws.Range("K134").Value = 36558.999
$ws.Range("M134").Value = -31488.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 16000
$ws.Range("J5").Value = 16000
$ws.Range("L5").Value = 16000
$ws.Range("N5").Value = -16224
$ws.Range("H10").Value = 337714.66
$ws.Range("I10").Value = 752637.25
$ws.Range("J10").Value = 5776.6
$ws.Range("K10").Value = 752637.25
$ws.Range("L10").Value = 5776.6
$ws.Range("M10").Value = -752468.25
$ws.Range("N10").Value = -6114.6
$ws.Range("H70").Value = 10725.533
$ws.Range("I70").Value = 10657.167
$ws.Range("K70").Value = 10657.167
$ws.Range("M70").Value = -10387.167
$ws.Range("H73").Value = 10725.533
$ws.Range("I73").Value = 10657.167
$ws.Range("K73").Value = 10657.167
$ws.Range("M73").Value = -9721.166999999999
$ws.Range("H80").Value = 5587.8887
$ws.Range("I80").Value = 3048.8333
$ws.Range("J80").Value = 10666
$ws.Range("K80").Value = 3048.8333
$ws.Range("L80").Value = 10666
$ws.Range("M80").Value = -2050.8333
$ws.Range("N80").Value = -12662
$ws.Range("H83").Value = 5587.8887
$ws.Range("I83").Value = 3048.8333
$ws.Range("J83").Value = 10666
$ws.Range("K83").Value = 15244.1665
$ws.Range("L83").Value = 53330
$ws.Range("M83").Value = -10252.1665
$ws.Range("N83").Value = -63314
$ws.Range("H126").Value = 5233
$ws.Range("I126").Value = 7199.3335
$ws.Range("K126").Value = 21598.0005
$ws.Range("M126").Value = -19128.0005
$ws.Range("H132").Value = 1734.2727
$ws.Range("I132").Value = 1314.3889
$ws.Range("K132").Value = 3943.1667
$ws.Range("M132").Value = -1413.1667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5429.2
$ws.Range("I61").Value = 5245.5713
$ws.Range("K61").Value = 5245.5713
$ws.Range("M61").Value = -5043.5713
$ws.Range("H68").Value = 1776.25
$ws.Range("I68").Value = 1307
$ws.Range("J68").Value = 2057.8
$ws.Range("K68").Value = 1307
$ws.Range("L68").Value = 2057.8
$ws.Range("M68").Value = -558
$ws.Range("N68").Value = -3555.8
$ws.Range("H71").Value = 1776.25
$ws.Range("I71").Value = 1307
$ws.Range("J71").Value = 2057.8
$ws.Range("K71").Value = 6535
$ws.Range("L71").Value = 10289
$ws.Range("M71").Value = -2791
$ws.Range("N71").Value = -17777
$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52246
$ws.Range("H88").Value = 10085.5
$ws.Range("I88").Value = 10085.5
$ws.Range("K88").Value = 10085.5
$ws.Range("M88").Value = -9657.5
$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -161232
$ws.Range("H91").Value = 10085.5
$ws.Range("I91").Value = 10085.5
$ws.Range("K91").Value = 10085.5
$ws.Range("M91").Value = -8603.5
$ws.Range("H93").Value = 2942.8572
$ws.Range("I93").Value = 2942.8572
$ws.Range("K93").Value = 2942.8572
$ws.Range("M93").Value = -1694.8572
$ws.Range("H113").Value = 5429.2
$ws.Range("I113").Value = 5245.5713
$ws.Range("K113").Value = 5245.5713
$ws.Range("M113").Value = -3075.5713
$ws.Range("H132").Value = 5035
$ws.Range("I132").Value = 5122.6665
$ws.Range("K132").Value = 15367.9995
$ws.Range("M132").Value = -12837.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 2383.3333
$ws.Range("I23").Value = 3000
$ws.Range("J23").Value = 1150
$ws.Range("K23").Value = 3000
$ws.Range("L23").Value = 1150
$ws.Range("M23").Value = -2771
$ws.Range("N23").Value = -1608
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H75").Value = 27368.143
$ws.Range("J75").Value = 27368.143
$ws.Range("L75").Value = 27368.143
$ws.Range("N75").Value = -29240.143
$ws.Range("H78").Value = 27368.143
$ws.Range("J78").Value = 27368.143
$ws.Range("L78").Value = 82104.429
$ws.Range("N78").Value = -91464.429
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800
$ws.Range("H126").Value = 2677.3076
$ws.Range("I126").Value = 2100.1428
$ws.Range("K126").Value = 6300.428400000001
$ws.Range("M126").Value = -3830.428400000001
$ws.Range("H132").Value = 820.1
$ws.Range("I132").Value = 298.5
$ws.Range("K132").Value = 895.5
$ws.Range("M132").Value = 1634.5
$ws.Range("H136").Value = 1747.5
$ws.Range("I136").Value = 1510.8572
$ws.Range("J136").Value = 2299.6667
$ws.Range("K136").Value = 4532.571599999999
$ws.Range("L136").Value = 6899.000100000001
$ws.Range("M136").Value = -1982.571599999999
$ws.Range("N136").Value = -11999.0001
